$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) values on sheet "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 203
    $ws.Range("F5").Value = 3514
    $ws.Range("F6").Value = 359
}
